$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.927.75"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.360.63"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.673"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "235.99"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.34"
$ws.Range("E7").Value = "  +11.45%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  +24.08%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.13"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.706.18"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.79"
$ws.Range("E14").Value = "  +8.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.79"
$ws.Range("E15").Value = "  +9.90%  "
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.361.97"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.916.99"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.23"
$ws.Range("E20").Value = "  +5.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.24"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.76"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("E26").Value = "  +7.29%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.79"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.58"
$ws.Range("E30").Value = "  +9.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.131"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +5.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.21"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  +4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.78"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.45"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0271"
$ws.Range("E39").Value = "  +7.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.24"
$ws.Range("E40").Value = "  +7.16%  "
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.186"
$ws.Range("E44").Value = "  +14.10%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0976"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.17"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.438.04"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.584.86"
$ws.Range("E51").Value = "  +0.49%  "
